# "Cartas funcionando y test parser ampliado, correr con bd vacia"
#
# Fix Pedro's (row 5) record: use the properly-accented surname and the
# canonically formatted street address, matching the formatting already
# used for the other rows in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5: "Perez Garcia" -> "Pérez García"
$ws.Range("B5").Value = "Pérez García"

# E5: "c/La playa 7" -> "C/ La playa 7"
$ws.Range("E5").Value = "C/ La playa 7"

# Leave the view scrolled/selected the way it was left after editing:
# column C at the left edge of the viewport, G5 selected.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G5").Select()
